$d = $word.ActiveDocument

# Remove the "Abstract Title" (AbstractTitle) custom paragraph style entirely.
$d.Styles("AbstractTitle").Delete()

# Update the "Abstract" style's paragraph spacing before from 5pt (100 twips)
# to 15pt (300 twips); SpaceBefore is expressed in points.
$d.Styles("Abstract").ParagraphFormat.SpaceBefore = 15

# Remove the "Footnote Block Text" (FootnoteBlockText) paragraph style entirely.
$d.Styles("FootnoteBlockText").Delete()
